# Update the "Handback" timestamps recorded in the handback-status report.
# These cells hold text values (formatted like "yyyy-mm-dd HH:mm:ss") rather
# than real date/time values, so we assign plain strings.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the
# 0cf647af-5bcc-4bb9-b199-699087b17e28.md row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-02 08:52:28"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the 0cf647af-5bcc-4bb9-b199-699087b17e28 row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-02 08:52:24"
$wsZhCn.Range("K3").Value = "2016-09-02 08:52:43"

# de-de sheet: Correspond Handback DateTime
# for the 0cf647af-5bcc-4bb9-b199-699087b17e28 row
# (its Correspond Handoff Datetime shares the same text as the Overview
# sheet's "Latest HO Xliff Generate Date" and is updated above)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-09-02 08:52:51"
